# Apply timetable updates to Section_A, Section_B, and Elective_Coordination sheets
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Section_A sheet
# ---------------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("B2").Value = "Free"
$wsA.Range("C2").Value = "CS151 (Elective)"
$wsA.Range("D2").Value = "CS151 (Elective)"
$wsA.Range("E2").Value = "C202"
$wsA.Range("F2").Value = "C202"

$wsA.Range("B3").Value = "EC161"
$wsA.Range("C3").Value = "Free"
$wsA.Range("D3").Value = "MA162"
$wsA.Range("E3").Value = "DS161"
$wsA.Range("F3").Value = "CS161"

$wsA.Range("B5").Value = "MA161"
$wsA.Range("C5").Value = "CS161"
$wsA.Range("D5").Value = "MA161"
$wsA.Range("E5").Value = "Free"
$wsA.Range("F5").Value = "DS161"

$wsA.Range("F6").Value = "CS151 (Tutorial)"

$wsA.Range("C7").Value = "C202"
$wsA.Range("E7").Value = "MA162"
$wsA.Range("F7").Value = "EC161"

$wsA.Range("B8").Value = "Free"

# ---------------------------------------------------------------------------
# Section_B sheet
# ---------------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("B2").Value = "C202"
$wsB.Range("C2").Value = "CS151 (Elective)"
$wsB.Range("D2").Value = "CS151 (Elective)"
$wsB.Range("E2").Value = "C202"
$wsB.Range("F2").Value = "EC161"

$wsB.Range("B3").Value = "CS161"
$wsB.Range("C3").Value = "CS161"
$wsB.Range("D3").Value = "MA162"
$wsB.Range("E3").Value = "MA161"
$wsB.Range("F3").Value = "DS161"

$wsB.Range("B5").Value = "MA162"
$wsB.Range("C5").Value = "EC161"
$wsB.Range("D5").Value = "C202"
$wsB.Range("E5").Value = "Free"

$wsB.Range("F6").Value = "CS151 (Tutorial)"

$wsB.Range("B7").Value = "EC161"
$wsB.Range("C7").Value = "DS161"
$wsB.Range("D7").Value = "Free"
$wsB.Range("F7").Value = "CS161"

$wsB.Range("B8").Value = "Free"

# ---------------------------------------------------------------------------
# Elective_Coordination sheet
# ---------------------------------------------------------------------------
$wsE = $wb.Worksheets.Item("Elective_Coordination")

$wsE.Range("C10").Value = "Wed"

$wsE.Range("C11").Value = "Tue"
$wsE.Range("D11").Value = "09:00-10:30"

$wsE.Range("C12").Value = "Fri"
$wsE.Range("D12").Value = "14:30-15:30"
